$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on all touched cells so numeric-looking strings
# (e.g. "1.00", "0.166") are preserved exactly as text, matching the
# source data which stores these as inline/shared strings, not numbers.

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '68.436.19'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -1.41%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '3.852.14'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -0.82%  '
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.03%  '
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.33%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '168.63'
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.00%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '3.852.16'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -0.77%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -1.05%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.166'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -2.08%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '6.48'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +1.17%  '
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -2.22%  '
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +4.28%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '37.06'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -3.13%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '4.497.94'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -0.90%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.855.98'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '68.534.86'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -1.36%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '18.50'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -1.40%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '7.37'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -3.10%  '
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -1.04%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '11.23'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +1.75%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '470.65'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -3.76%  '
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -1.44%  '
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  -3.51%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '83.44'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -2.18%  '
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -2.71%  '
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -1.95%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '10.18'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.16%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.95'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -1.15%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.001.83'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -0.82%  '
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -1.89%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '31.44'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -1.32%  '
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -3.65%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '9.34'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -3.03%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '3.817.81'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '3.76'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +10.31%  '
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -2.38%  '
$c = $ws.Range('B39')
$c.NumberFormat = '@'
$c.Value = 'Mantle'
$c = $ws.Range('C39')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '1.02'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -1.07%  '
$c = $ws.Range('B40')
$c.NumberFormat = '@'
$c.Value = 'Kaspa'
$c = $ws.Range('C40')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.140'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -2.00%  '
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -2.61%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -3.66%  '
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -4.50%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '8.70'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +0.13%  '
$c = $ws.Range('B46')
$c.NumberFormat = '@'
$c.Value = 'FLOKI'
$c = $ws.Range('C46')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.000295'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +6.87%  '
$c = $ws.Range('B47')
$c.NumberFormat = '@'
$c.Value = 'Bittensor'
$c = $ws.Range('C47')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '418.03'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -3.98%  '
$c = $ws.Range('B48')
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c = $ws.Range('C48')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '46.97'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -2.07%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '141.75'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +0.43%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '26.08'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +3.67%  '
